$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Update values: B3 changes from 80000 to 47000
$ws.Range("B3").Value = 47000

# B4 changes from a number (80000) to a text value of two spaces,
# which breaks the downstream formulas (B7, B10, B11) into #VALUE! errors.
$ws.Range("B4").Value = "  "

# Update the active selection / view (drop the frozen topLeftCell scroll,
# move selection to B4).
$ws.Activate()
$ws.Range("B4").Select()
